# hx_vcenter_register_input.xlsx — add prompts for username/password for HX and vCenter
# (per the commit message the sheet is actually trimmed down to just the two
#  FQDN/IP prompt columns: HyperFlex_FQDN/IP and vCenter_FQDN/IP)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the vCenter_FQDN/IP column (D) next to column A, preserving its
# original column width exactly (Cut/Insert keeps the col metadata intact).
$ws.Columns("D").Cut() | Out-Null
$ws.Columns("B").Insert() | Out-Null

# Drop every other column (old B/C HyperFlex user+password, and the old
# vCenter user/password/datacenter/cluster columns E:H).
$ws.Columns("C:H").Delete() | Out-Null

# The vCenter_User hyperlink (mailto:) no longer exists in the sheet.
$ws.Hyperlinks.Delete() | Out-Null

# Drop the now-unused "Hyperlink" cell style.
$wb.Styles("Hyperlink").Delete() | Out-Null

# Match the saved selection state (full-column selection starting at C1).
$ws.Columns("C:F").Select() | Out-Null
